$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text that happens to
# look numeric/percentage-like (e.g. "270.20", "3.44%"); Excel's automatic
# type inference would otherwise coerce these into real numbers/percentages
# and silently drop significant trailing zeros. Force each target cell to
# Text format (single-area ranges only -- multi-area NumberFormat writes are
# unreliable) before writing the literal string value.

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "270.20"
$ws.Range("E2").Value = "3.44%"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "26.75"
$ws.Range("E3").Value = "-1.39%"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "4.724"
$ws.Range("E4").Value = "0.25%"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06120"
$ws.Range("E5").Value = "-1.39%"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "6.745"
$ws.Range("E6").Value = "0.30%"

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8568"
$ws.Range("E7").Value = "0.76%"

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8931"
$ws.Range("E8").Value = "-1.54%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.04%"

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05051"
$ws.Range("E10").Value = "7.94%"

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07142"
$ws.Range("E11").Value = "0.81%"

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03170"
$ws.Range("E12").Value = "0.08%"

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09043"
$ws.Range("E13").Value = "-0.20%"

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001546"
$ws.Range("E14").Value = "0.41%"

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006077"
$ws.Range("E15").Value = "-1.35%"

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006107"
$ws.Range("E16").Value = "-0.32%"

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").Value = "-0.06%"

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.176"
$ws.Range("E18").Value = "0.18%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.97%"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1281"
$ws.Range("E21").Value = "-1.42%"

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "3.846"
$ws.Range("E22").Value = "-6.42%"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04243"
$ws.Range("E23").Value = "0.39%"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").Value = "-3.29%"

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004152"
$ws.Range("E25").Value = "0.43%"

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.01%"

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001682"
$ws.Range("E27").Value = "4.05%"

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03962"
$ws.Range("E40").Value = "1.11%"

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1120"
$ws.Range("E41").Value = "0.54%"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004195"
$ws.Range("E42").Value = "1.48%"

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002039"
$ws.Range("E43").Value = "-6.65%"

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01175"
$ws.Range("E44").Value = "-13.09%"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005129"
$ws.Range("E45").Value = "-0.87%"

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.01%"

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9202"
$ws.Range("E47").Value = "449.42%"

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02994"
$ws.Range("E48").Value = "-16.63%"

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "0.01%"

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").Value = "0.01%"
